{"js": "// Bold + color (\"#2C3E50\") the quantitative metrics (percentages, dollar\n// amounts, etc.) inside specific bullet lines, matching the commit's\n// \"hybrid bold + color highlighting for impact metrics\" change.\n//\n// Strategy: locate each target paragraph by a unique substring (so we never\n// touch the other paragraphs in the doc that happen to repeat the same\n// numbers, e.g. \"73.5%\" / \"$4.7M\" / \"23% to 64%\" also appear in the summary\n// / key-projects sections and must stay untouched), then for every metric\n// substring inside that paragraph call `paragraph.search(metric)` \u2014 scoped\n// to just that paragraph, so it can't match text elsewhere in the document \u2014\n// and set Font.bold / Font.color on the single hit. Word automatically\n// splits the run the same way: plain runs around bolded/colored metric runs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\n// Map of \"unique substring that identifies the paragraph\" -> \"ordered list\n// of metric substrings inside it to bold + color\".\nconst targets = [\n  [\"Discovered systematic race\", [\"23%\", \"64%\"]],\n  [\"Utilized advanced sampling\", [\"\\u00B14.2%\", \"\\u00B12.1%\", \"71%\", \"87%\"]],\n  [\"Trigonometric algorithm\", [\"73.5%\", \"$4.7M\"]],\n  [\"Built real-time FEC\", [\"$2\"]],\n  [\"Algorithmic innovation\", [\"73.5%\"]],\n  [\"$4.7M savings enabled\", [\"$4.7M\"]],\n  [\"178% accuracy improvement\", [\"178%\"]],\n];\n\nfor (const p of paragraphs.items) {\n  for (const [key, metrics] of targets) {\n    if (p.text.includes(key)) {\n      for (const metric of metrics) {\n        const found = p.search(metric, { matchCase: true });\n        found.load(\"items\");\n        await context.sync();\n        if (found.items.length > 0) {\n          const hit = found.items[0];\n          hit.font.bold = true;\n          hit.font.color = \"#2C3E50\";\n        }\n      }\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Bold + color (\"2C3E50\") the quantitative metrics (percentages, dollar\n# amounts, etc.) inside specific bullet lines, matching the commit's\n# \"hybrid bold + color highlighting for impact metrics\" change.\n#\n# Strategy: for each target paragraph (located by a unique substring so we\n# never touch the other paragraphs in the doc that happen to repeat the\n# same numbers, e.g. \"73.5%\" / \"$4.7M\" / \"23% to 64%\" appear elsewhere too),\n# walk the paragraph's Range left-to-right, Find.Execute each metric in\n# sequence (narrowing the search window to [cursor, paragraphEnd] so we\n# never re-match an earlier occurrence), then set Font.Bold / Font.Color\n# on just that found sub-range. This naturally splits the run exactly the\n# way Word would: plain text runs around bolded/colored metric runs.\n\nfunction Format-MetricsInParagraph($doc, $para, $metrics) {\n    $pStart = $para.Range.Start\n    $pEnd = $para.Range.End\n    $cursor = $pStart\n    foreach ($metric in $metrics) {\n        $searchRange = $doc.Range($cursor, $pEnd)\n        $searchRange.Find.MatchWildcards = $false\n        $found = $searchRange.Find.Execute($metric)\n        if ($found) {\n            $searchRange.Font.Bold = 1\n            $searchRange.Font.Color = \"2C3E50\"\n            $cursor = $searchRange.End\n        }\n    }\n}\n\n$d = $word.ActiveDocument\n\n# Map of \"unique substring that identifies the paragraph\" -> \"ordered list\n# of metric substrings inside it to bold+color\".\n$targets = @(\n    @(\"\u2022 Discovered systematic race\", @(\"23%\", \"64%\")),\n    @(\"\u2022 Utilized advanced sampling\", @([char]0x00B1 + \"4.2%\", [char]0x00B1 + \"2.1%\", \"71%\", \"87%\")),\n    @(\"\u2022 Trigonometric algorithm\", @(\"73.5%\", \"$4.7M\")),\n    @(\"\u2022 Built real-time FEC\", @(\"$2\")),\n    @(\"\u2022 Algorithmic innovation\", @(\"73.5%\")),\n    @(\"\u2022 $4.7M savings enabled\", @(\"$4.7M\")),\n    @(\"\u2022 178% accuracy improvement\", @(\"178%\"))\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    foreach ($target in $targets) {\n        $key = $target[0]\n        $metrics = $target[1]\n        if ($text -like \"*$key*\") {\n            Format-MetricsInParagraph $d $p $metrics\n        }\n    }\n}\n"}
